# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" table (rows 16-98, sheet Hoja1) was
# re-sorted from most-recent-first (2502 .. 1804) to chronological,
# oldest-first order (1804 .. 2502), and the "Valor Mora" amounts were
# updated to match the new first/last period's installment values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Build the ordered list of periods 1804 .. 2502 (YYMM, ascending).
$periods = @()
$ranges = @(
    @(18, 4, 12),
    @(19, 1, 12),
    @(20, 1, 12),
    @(21, 1, 12),
    @(22, 1, 12),
    @(23, 1, 12),
    @(24, 1, 12),
    @(25, 1, 2)
)
foreach ($r in $ranges) {
    $yy = $r[0]; $mStart = $r[1]; $mEnd = $r[2]
    for ($m = $mStart; $m -le $mEnd; $m++) {
        $periods += ("{0:d2}{1:d2}" -f $yy, $m)
    }
}

$firstRow = 16
$lastRow = 98

# Valor Mora per period: the first period (1804) keeps the old last row's
# value, the next 13 periods (1805-1905) keep the old "13" segment value,
# the middle block (1906-2501) keeps the steady installment value, and the
# final period (2502) keeps the old first row's value.
for ($i = 0; $i -lt $periods.Count; $i++) {
    $row = $firstRow + $i
    $period = $periods[$i]

    if ($i -eq 0) {
        $valorMora = 5208
    } elseif ($i -ge 1 -and $i -le 13) {
        $valorMora = 7812
    } elseif ($i -eq ($periods.Count - 1)) {
        $valorMora = 19791
    } else {
        $valorMora = 31249
    }

    $ws.Range("E$row").Value = $period
    $ws.Range("F$row").Value = $valorMora
}
